# Update the FRCcombined (column W) measurements on the "Sheet1" tab.
# These values were recomputed after data files were renamed / a new
# data file was added upstream; the rest of the workbook is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("W2").Value = 8.097165991902834
$ws.Range("W4").Value = 8.5106382978723403
$ws.Range("W5").Value = 6.7340067340067336
$ws.Range("W6").Value = 7.7821011673151741
$ws.Range("W7").Value = 7.5757575757575761
$ws.Range("W8").Value = 8.2987551867219906
$ws.Range("W9").Value = 7.4074074074074066
$ws.Range("W10").Value = 7.4074074074074066
$ws.Range("W11").Value = 10.256410256410255
$ws.Range("W12").Value = 8.695652173913043
$ws.Range("W13").Value = 8.695652173913043
$ws.Range("W14").Value = 10.1010101010101
$ws.Range("W15").Value = 9.661835748792269
$ws.Range("W16").Value = 8.4388185654008439
$ws.Range("W17").Value = 7.4626865671641784
$ws.Range("W18").Value = 7.8740157480314954
$ws.Range("W19").Value = 9.5693779904306204
